{"js": "const replacements = [\n  [\"2024-05-16 Thursday\", \"2024-05-17 Friday\"],\n  [\"125\u00f74=\", \"827\u00f73=\"],\n  [\"854\u00f76=\", \"402\u00f73=\"],\n  [\"542\u00f75=\", \"152\u00f76=\"],\n  [\"755\u00f77=\", \"376\u00f78=\"],\n  [\"566\u00f75=\", \"974\u00f75=\"],\n  [\"310\u00f77=\", \"515\u00f74=\"],\n  [\"701\u00f75=\", \"262\u00f78=\"],\n  [\"794\u00f79=\", \"878\u00f78=\"],\n  [\"299\u00f72=\", \"891\u00f77=\"],\n  [\"826\u00f72=\", \"790\u00f75=\"],\n  [\"126\u00f75=\", \"132\u00f77=\"],\n  [\"952\u00f79=\", \"243\u00f73=\"],\n  [\"396\u00f77=\", \"433\u00f73=\"],\n  [\"390\u00f72=\", \"198\u00f76=\"],\n  [\"340\u00f73=\", \"221\u00f75=\"],\n  [\"925\u00f73=\", \"921\u00f78=\"],\n  [\"157\u00f72=\", \"844\u00f75=\"],\n  [\"396\u00f73=\", \"516\u00f74=\"],\n  [\"559\u00f77=\", \"184\u00f74=\"],\n  [\"734\u00f76=\", \"915\u00f76=\"],\n  [\"237\u00f72=\", \"991\u00f75=\"],\n  [\"469\u00f72=\", \"191\u00f74=\"],\n  [\"483\u00f74=\", \"316\u00f78=\"],\n  [\"912\u00f76=\", \"878\u00f78=\"],\n  [\"223\u00f79=\", \"305\u00f79=\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + before);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"2024-05-16 Thursday\", \"2024-05-17 Friday\")\n  ,@(\"125\u00f74=\", \"827\u00f73=\")\n  ,@(\"854\u00f76=\", \"402\u00f73=\")\n  ,@(\"542\u00f75=\", \"152\u00f76=\")\n  ,@(\"755\u00f77=\", \"376\u00f78=\")\n  ,@(\"566\u00f75=\", \"974\u00f75=\")\n  ,@(\"310\u00f77=\", \"515\u00f74=\")\n  ,@(\"701\u00f75=\", \"262\u00f78=\")\n  ,@(\"794\u00f79=\", \"878\u00f78=\")\n  ,@(\"299\u00f72=\", \"891\u00f77=\")\n  ,@(\"826\u00f72=\", \"790\u00f75=\")\n  ,@(\"126\u00f75=\", \"132\u00f77=\")\n  ,@(\"952\u00f79=\", \"243\u00f73=\")\n  ,@(\"396\u00f77=\", \"433\u00f73=\")\n  ,@(\"390\u00f72=\", \"198\u00f76=\")\n  ,@(\"340\u00f73=\", \"221\u00f75=\")\n  ,@(\"925\u00f73=\", \"921\u00f78=\")\n  ,@(\"157\u00f72=\", \"844\u00f75=\")\n  ,@(\"396\u00f73=\", \"516\u00f74=\")\n  ,@(\"559\u00f77=\", \"184\u00f74=\")\n  ,@(\"734\u00f76=\", \"915\u00f76=\")\n  ,@(\"237\u00f72=\", \"991\u00f75=\")\n  ,@(\"469\u00f72=\", \"191\u00f74=\")\n  ,@(\"483\u00f74=\", \"316\u00f78=\")\n  ,@(\"912\u00f76=\", \"878\u00f78=\")\n  ,@(\"223\u00f79=\", \"305\u00f79=\")\n)\n\nforeach ($pair in $replacements) {\n  $findText = $pair[0]\n  $replaceText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $findText\n  $find.Replacement.Text = $replaceText\n  $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}"}
